# Updates the cryptos price-table cells per the target diff.
#
# Several "Price" (column D) values are plain numeric-looking strings
# (e.g. "306.55") that Excel would otherwise auto-convert to a Number
# the moment `.Value` is assigned. The source workbook stores every
# D/E cell as literal TEXT (inline string), so for every cell whose new
# text *could* parse as a number we briefly force a text number format,
# assign the value, then clear the formatting again (this reverts the
# cell's style index back to the sheet default, since the original
# cells carry no explicit style).
#
# Cells whose new text can never parse as a plain number (multi-dot
# "thousands" prices like "26.442.11", the padded percentage strings in
# column E, coin names, links, ...) are just assigned directly - Excel
# already keeps those as text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, [string]$text)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

function Set-Cell {
    param([int]$row, [int]$col, [string]$text)
    $cell = $ws.Cells.Item($row, $col)
    Set-TextValue $cell $text
}

Set-Cell 2 4 "26.442.11"
Set-Cell 2 5 "  -3.81%  "
Set-Cell 3 4 "1.770.85"
Set-Cell 3 5 "  -3.01%  "
Set-Cell 4 5 "  +0.10%  "
Set-Cell 5 5 "  +0.10%  "
Set-Cell 6 4 "306.55"
Set-Cell 6 5 "  -1.97%  "
Set-Cell 7 4 "0.4305"
Set-Cell 7 5 "  +1.18%  "
Set-Cell 8 4 "0.3663"
Set-Cell 8 5 "  +1.39%  "
Set-Cell 9 4 "0.07230"
Set-Cell 9 5 "  +0.36%  "
Set-Cell 10 4 "0.8503"
Set-Cell 10 5 "  -1.45%  "
Set-Cell 11 4 "20.33"
Set-Cell 11 5 "  -1.43%  "
Set-Cell 12 4 "1.786.50"
Set-Cell 12 5 "  -2.08%  "
Set-Cell 13 4 "6.437"
Set-Cell 13 5 "  -0.61%  "
Set-Cell 14 4 "5.240"
Set-Cell 14 5 "  -2.76%  "
Set-Cell 15 4 "0.06847"
Set-Cell 15 5 "  -1.01%  "
Set-Cell 16 5 "  +0.26%  "
Set-Cell 17 5 "  -1.56%  "
Set-Cell 18 4 "0.000008692"
Set-Cell 18 5 "  -2.82%  "
Set-Cell 19 5 "  +0.15%  "
Set-Cell 20 4 "15.06"
Set-Cell 20 5 "  -2.12%  "
Set-Cell 21 4 "26.439.81"
Set-Cell 21 5 "  -3.81%  "
Set-Cell 22 4 "5.101"
Set-Cell 22 5 "  -0.40%  "
Set-Cell 23 4 "11.29"
Set-Cell 23 5 "  +3.83%  "
Set-Cell 24 4 "1.994.92"
Set-Cell 24 5 "  -3.38%  "
Set-Cell 25 4 "152.13"
Set-Cell 25 5 "  -2.10%  "
Set-Cell 26 4 "1.849"
Set-Cell 26 5 "  -6.92%  "
Set-Cell 27 5 "  -2.89%  "
Set-Cell 28 4 "5.093"
Set-Cell 28 5 "  -1.21%  "
Set-Cell 29 4 "114.73"
Set-Cell 29 5 "  +0.40%  "
Set-Cell 30 4 "1.722"
Set-Cell 30 5 "  -3.88%  "
Set-Cell 31 4 "0.08970"
Set-Cell 31 5 "  +1.07%  "
Set-Cell 32 4 "0.7246"
Set-Cell 32 5 "  -3.19%  "
Set-Cell 33 4 "1.117"
Set-Cell 33 5 "  -0.43%  "
Set-Cell 34 4 "4.330"
Set-Cell 35 2 "Frax"
Set-Cell 35 3 "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-Cell 35 4 "1.002"
Set-Cell 35 5 "  +0.10%  "
Set-Cell 36 2 "HuobiToken"
Set-Cell 36 3 "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-Cell 36 4 "2.745"
Set-Cell 36 5 "  -7.61%  "
Set-Cell 37 4 "1.079"
Set-Cell 37 5 "  -0.56%  "
Set-Cell 38 4 "0.05168"
Set-Cell 38 5 "  -2.03%  "
Set-Cell 39 5 "  -1.41%  "
Set-Cell 40 5 "  -2.87%  "
Set-Cell 41 5 "  -3.14%  "
Set-Cell 42 5 "  -9.25%  "
Set-Cell 43 4 "6.242"
Set-Cell 43 5 "  -3.03%  "
Set-Cell 44 4 "8.030"
Set-Cell 44 5 "  -3.87%  "
Set-Cell 45 4 "104.94"
Set-Cell 45 5 "  -1.30%  "
Set-Cell 46 2 "PaxDollar"
Set-Cell 46 3 "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-Cell 46 4 "1.002"
Set-Cell 46 5 "  +0.16%  "
Set-Cell 47 2 "EnergySwap"
Set-Cell 47 3 "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-Cell 47 4 "10.15"
Set-Cell 47 5 "  -3.24%  "
Set-Cell 48 5 "  -4.10%  "
Set-Cell 49 4 "0.4481"
Set-Cell 49 5 "  -4.15%  "
Set-Cell 50 4 "1.587"
Set-Cell 50 5 "  -1.71%  "
Set-Cell 51 4 "1.743"
Set-Cell 51 5 "  +2.24%  "
